$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wb.Worksheets.Item("Listado de premios").Name = "Ganadores y perdedores"
$wb.Worksheets.Item("Problema 3").Name = "Problema 3 %"

$wsGanadores = $wb.Worksheets.Item("Ganadores y perdedores")
$wsAutos     = $wb.Worksheets.Item("Modelo autos 2025")
$wsProblema  = $wb.Worksheets.Item("Problema 3 %")

# --- "Modelo autos 2025": turn the repeated per-row formulas in columns B and C
#     into proper fill-down (shared) formulas across B2:B8 / C2:C8 ---
$wsAutos.Range("B2:B8").Formula = '=IF(A2="Mercedes 321","15060","7230")'
$wsAutos.Range("C2:C8").Formula = '=IF(B2="15060","Aplazado","Al contado")'

# --- Update each sheet's selection ---
$wsGanadores.Range("A1:G18").Select()
$wsAutos.Range("A1:E8").Select()
$wsProblema.Range("F26").Select()

# --- Make "Problema 3 %" the active sheet/tab ---
$wsProblema.Activate()
